$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2:B28 with the retrained 15-minute values
$ws.Cells.Item(2, 1).Value = 5142
$ws.Cells.Item(2, 2).Value = 45923
$ws.Cells.Item(3, 1).Value = 5136
$ws.Cells.Item(3, 2).Value = 45923.01041666666
$ws.Cells.Item(4, 1).Value = 5093
$ws.Cells.Item(4, 2).Value = 45923.02083333334
$ws.Cells.Item(5, 1).Value = 5035
$ws.Cells.Item(5, 2).Value = 45923.03125
$ws.Cells.Item(6, 1).Value = 5049
$ws.Cells.Item(6, 2).Value = 45923.04166666666
$ws.Cells.Item(7, 1).Value = 5056
$ws.Cells.Item(7, 2).Value = 45923.05208333334
$ws.Cells.Item(8, 1).Value = 5058
$ws.Cells.Item(8, 2).Value = 45923.0625
$ws.Cells.Item(9, 1).Value = 5015
$ws.Cells.Item(9, 2).Value = 45923.07291666666
$ws.Cells.Item(10, 1).Value = 5028
$ws.Cells.Item(10, 2).Value = 45923.08333333334
$ws.Cells.Item(11, 1).Value = 4988
$ws.Cells.Item(11, 2).Value = 45923.09375
$ws.Cells.Item(12, 1).Value = 4992
$ws.Cells.Item(12, 2).Value = 45923.10416666666
$ws.Cells.Item(13, 1).Value = 4986
$ws.Cells.Item(13, 2).Value = 45923.11458333334
$ws.Cells.Item(14, 1).Value = 5040
$ws.Cells.Item(14, 2).Value = 45923.125
$ws.Cells.Item(15, 1).Value = 5053
$ws.Cells.Item(15, 2).Value = 45923.13541666666
$ws.Cells.Item(16, 1).Value = 5109
$ws.Cells.Item(16, 2).Value = 45923.14583333334
$ws.Cells.Item(17, 1).Value = 5027
$ws.Cells.Item(17, 2).Value = 45923.15625
$ws.Cells.Item(18, 1).Value = 5128
$ws.Cells.Item(18, 2).Value = 45923.16666666666
$ws.Cells.Item(19, 1).Value = 5163
$ws.Cells.Item(19, 2).Value = 45923.17708333334
$ws.Cells.Item(20, 1).Value = 5278
$ws.Cells.Item(20, 2).Value = 45923.1875
$ws.Cells.Item(21, 1).Value = 5314
$ws.Cells.Item(21, 2).Value = 45923.19791666666
$ws.Cells.Item(22, 1).Value = 5556
$ws.Cells.Item(22, 2).Value = 45923.20833333334
$ws.Cells.Item(23, 1).Value = 5692
$ws.Cells.Item(23, 2).Value = 45923.21875
$ws.Cells.Item(24, 1).Value = 5780
$ws.Cells.Item(24, 2).Value = 45923.22916666666
$ws.Cells.Item(25, 1).Value = 5832
$ws.Cells.Item(25, 2).Value = 45923.23958333334
$ws.Cells.Item(26, 1).Value = 6079
$ws.Cells.Item(26, 2).Value = 45923.25
$ws.Cells.Item(27, 1).Value = 6252
$ws.Cells.Item(27, 2).Value = 45923.26041666666
$ws.Cells.Item(28, 1).Value = 6295
$ws.Cells.Item(28, 2).Value = 45923.27083333334

# Remove the now-unused trailing rows (29:42) so the sheet shrinks to A1:B28
$ws.Range("A29:B42").EntireRow.Delete()
